$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected in the source file; unprotect to allow cell edits, restore protection at the end.
$ws.Unprotect()

# Update the confidential disclaimer string (A37): date 2021-03-24 -> 2021-03-25
$ws.Range("A37").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-25 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for holdings rows 2-34
$ws.Range("D2").Value2 = 0.03875852259849078
$ws.Range("E2").Value2 = -0.001557026080186819
$ws.Range("D3").Value2 = 0.02191955458746976
$ws.Range("E3").Value2 = 0.001889768698558392
$ws.Range("D4").Value2 = 0.02013868640143529
$ws.Range("E4").Value2 = -0.0004140786749481373
$ws.Range("D5").Value2 = 0.0411233424670041
$ws.Range("E5").Value2 = -0.003468609087755703
$ws.Range("D6").Value2 = 0.03769151455713327
$ws.Range("E6").Value2 = -0.0005873140172280022
$ws.Range("D7").Value2 = 0.0211533562417822
$ws.Range("E7").Value2 = 0.000583544057576546
$ws.Range("D8").Value2 = 0.03761986498635776
$ws.Range("E8").Value2 = 0.004294755877034095
$ws.Range("D9").Value2 = 0.0216028393454298
$ws.Range("E9").Value2 = -0.00146279027244467
$ws.Range("D10").Value2 = 0.02589204818338614
$ws.Range("E10").Value2 = 0.01207755058798599
$ws.Range("D11").Value2 = 0.02387176070726961
$ws.Range("E11").Value2 = 0.008905486929043471
$ws.Range("D12").Value2 = 0.05784539768988924
$ws.Range("E12").Value2 = 0.008356545961002881
$ws.Range("D13").Value2 = 0.02667273225087122
$ws.Range("E13").Value2 = 0.0003702332469455172
$ws.Range("D14").Value2 = 0.02730243651837606
$ws.Range("E14").Value2 = 0.01136925358378638
$ws.Range("D15").Value2 = 0.03502994696402215
$ws.Range("E15").Value2 = 0.006490089457990011
$ws.Range("D16").Value2 = 0.01895021423320413
$ws.Range("E16").Value2 = 0.001684919966301601
$ws.Range("D17").Value2 = 0.03006703465813068
$ws.Range("E17").Value2 = 0.004379162485174781
$ws.Range("D18").Value2 = 0.02423587877866093
$ws.Range("E18").Value2 = -0.000692680674209023
$ws.Range("D19").Value2 = 0.1345842276024523
$ws.Range("E19").Value2 = -0.001338688085676054
$ws.Range("D20").Value2 = 0.009620134406299665
$ws.Range("E20").Value2 = -0.002229800629590772
$ws.Range("D21").Value2 = 0.01592085282502484
$ws.Range("E21").Value2 = 0.01157826181344523
$ws.Range("D22").Value2 = 0.01724006077209064
$ws.Range("E22").Value2 = 0.01347038139031653
$ws.Range("D23").Value2 = 0.01686343419982886
$ws.Range("E23").Value2 = -0.007801418439716157
$ws.Range("D24").Value2 = 0.02133297392533429
$ws.Range("E24").Value2 = 0.003317483136127342
$ws.Range("D25").Value2 = 0.01169621637971994
$ws.Range("E25").Value2 = 0.007228158390949035
$ws.Range("D26").Value2 = 0.04313210895547678
$ws.Range("E26").Value2 = 0.01036002518459167
$ws.Range("D27").Value2 = 0.02572855988253544
$ws.Range("E27").Value2 = 0
$ws.Range("D28").Value2 = 0.04809781921551032
$ws.Range("E28").Value2 = -0.0009718172983479434
$ws.Range("D29").Value2 = 0.05672166249984729
$ws.Range("E29").Value2 = 0.004805687203791376
$ws.Range("D30").Value2 = 0.01334185231002792
$ws.Range("E30").Value2 = 0.01979522184300331
$ws.Range("D31").Value2 = 0.01458639328479244
$ws.Range("E31").Value2 = -0.0009402914903620108
$ws.Range("D32").Value2 = 0.04455659678486892
$ws.Range("E32").Value2 = 0.0005208333333333037
$ws.Range("D33").Value2 = 0.01670197578727731
$ws.Range("E33").Value2 = 0.005150490906164373
$ws.Range("E34").Value2 = 0.002924418377586946

# Restore sheet protection (best-effort; engine does not preserve the original legacy password hash/flags)
$ws.Protect()

Write-Host "Edit complete"
